# Refresh the cryptocurrency price/volume snapshot (Price = column D, Volume(1h) = column E).
# Both columns are stored as plain text in the sheet (see the original inlineStr cells), so
# every write below targets Range.Value with a string. Several "Price" figures (e.g. 516.29)
# are plain decimals that Excel would otherwise auto-convert to a Double; those are written
# with a leading apostrophe (the normal Excel "treat as text" prefix) so they stay text,
# exactly like the untouched text cells already on the sheet (e.g. "56.867.99").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell reference -> new value (apostrophe-prefixed where Excel would otherwise read it as a number)
$updates = [ordered]@{
    'D2' = '57.162.09'
    'E2' = '  -0.53%  '
    'D3' = '3.022.02'
    'E3' = '  +0.46%  '
    'E4' = '  -0.08%  '
    'D5' = '''516.29'
    'E5' = '  +1.45%  '
    'D6' = '''138.72'
    'E6' = '  -0.21%  '
    'D7' = '''0.997'
    'E7' = '  -0.20%  '
    'D8' = '''0.432'
    'E8' = '  -1.19%  '
    'D9' = '''7.30'
    'E9' = '  -2.76%  '
    'D10' = '''0.109'
    'E10' = '  -0.88%  '
    'D11' = '''0.369'
    'E11' = '  +1.27%  '
    'D12' = '3.513.48'
    'E12' = '  -0.23%  '
    'E13' = '  -3.54%  '
    'D14' = '''26.65'
    'E14' = '  +0.89%  '
    'E15' = '  +2.97%  '
    'D16' = '57.013.68'
    'E16' = '  -0.77%  '
    'D17' = '''6.18'
    'E17' = '  -0.17%  '
    'D18' = '3.015.34'
    'E18' = '  +0.13%  '
    'D19' = '''13.31'
    'E19' = '  +3.93%  '
    'D20' = '''8.10'
    'E20' = '  +1.87%  '
    'D21' = '''327.06'
    'E21' = '  -1.10%  '
    'E22' = '  +0.53%  '
    'D23' = '''0.503'
    'E23' = '  +1.09%  '
    'D24' = '''64.40'
    'E24' = '  -0.01%  '
    'D25' = '3.121.73'
    'E25' = '  -0.32%  '
    'D26' = '''0.999'
    'E26' = '  +0.12%  '
    'E27' = '  -3.63%  '
    'D28' = '0.0₃0886'
    'E28' = '  -3.76%  '
    'D29' = '''6.62'
    'E29' = '  -2.36%  '
    'D30' = '''7.23'
    'E30' = '  -1.45%  '
    'E31' = '  -0.06%  '
    'E32' = '  +2.54%  '
    'D33' = '''20.64'
    'E33' = '  +0.46%  '
    'D34' = '''153.79'
    'E34' = '  +0.17%  '
    'D35' = '''4.58'
    'E35' = '  -2.82%  '
    'D36' = '''5.82'
    'E36' = '  -0.82%  '
    'D37' = '''25.29'
    'E37' = '  +3.77%  '
    'D38' = '''1.26'
    'E38' = '  -1.42%  '
    'D39' = '''0.0671'
    'E39' = '  -1.68%  '
    'D40' = '''36.99'
    'E40' = '  -0.79%  '
    'D41' = '''3.86'
    'E41' = '  +0.65%  '
    'E42' = '  -0.10%  '
    'D43' = '''0.658'
    'E43' = '  +1.39%  '
    'E44' = '  -1.74%  '
    'D45' = '''6.16'
    'E45' = '  +2.37%  '
    'D46' = '2.172.25'
    'E46' = '  -3.78%  '
    'D47' = '''0.948'
    'E47' = '  -3.66%  '
    'D48' = '''0.0243'
    'E48' = '  +1.20%  '
    'D49' = '''19.64'
    'E49' = '  +1.40%  '
    'E50' = '  -4.53%  '
    'D51' = '''0.183'
    'E51' = '  +0.42%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
